$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Schedule")
$ws2 = $wb.Worksheets.Item("Schedule_date")

# ---------------------------------------------------------------------------
# Sheet "Schedule" (sheet1): add the combined "Manhattan plots" topic, the new
# ggplot-extensions / complexheatmap topic, the Thanksgiving no-class week,
# and fill in the previously-blank topic for the capstone-prep week.
# ---------------------------------------------------------------------------

# Week 8's "Open session, capstone prep" row was missing a Topic value.
$ws1.Range("C9").Value = "Open session, capstone prep"

# Week 11: "Manhattan plots" lecture grows to include the recitation content.
$ws1.Range("C12").Value = "Manhattan plots and making lots of plots at once"

# Week 13 used to be "Making lots of plots at once" -- that recitation slot
# now covers ggplot extension packages / ComplexHeatmap instead.
$ws1.Range("C14").Value = "ggplot extension packages and complexheatmap"

# Insert the new Thanksgiving week (week 14) before the old week-14 row,
# pushing the two capstone rows down to weeks 15 and 16.
$ws1.Rows.Item(15).EntireRow.Insert()
$ws1.Range("A15").Value = 14
$ws1.Range("B15").Value = "No class, Thanksgiving"
$ws1.Range("C15").Value = "Relaxing and eating"
$ws1.Range("A16").Value = 15
$ws1.Range("A17").Value = 16

# Column B needs to be widened to fit the longer module text.
$ws1.Columns.Item(2).ColumnWidth = 23.998697916666668

# ---------------------------------------------------------------------------
# Sheet "Schedule_date" (sheet2): same content changes, mirrored against the
# dated version of the schedule -- plus it was missing the whole
# "Open session, capstone prep" row that "Schedule" already had.
# ---------------------------------------------------------------------------

# Insert week 8's "Open session, capstone prep" row (present in "Schedule"
# but missing here), pushing everything after it down by one row.
$ws2.Rows.Item(9).EntireRow.Insert()
$ws2.Range("C9").Value = "Open session, capstone prep"
$ws2.Range("D9").Value = "Open session, capstone prep"

# Week 11 (now row 12): combined Manhattan-plots topic.
$ws2.Range("D12").Value = "Manhattan plots and making lots of plots at once"

# Week 13 (now row 14): new ggplot-extensions / complexheatmap topic.
$ws2.Range("D14").Value = "ggplot extension packages and complexheatmap"

# Insert the new Thanksgiving week (now row 15) before the old week-14 row.
$ws2.Rows.Item(15).EntireRow.Insert()
$ws2.Range("C15").Value = "No class, Thanksgiving"
$ws2.Range("D15").Value = "Relaxing and eating"

# Re-number the Week/Date columns for every data row now that the two new
# rows have shifted things around -- straight weekly sequence, 16 weeks.
$weeks = 1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16
$dates = 44796,44803,44810,44817,44825,44832,44838,44845,44852,44859,44866,44873,44880,44887,44894,44900
for ($i = 0; $i -lt 16; $i++) {
    $r = $i + 2
    $ws2.Range("A$r").Value = $weeks[$i]
    $ws2.Range("B$r").Value = $dates[$i]
    $ws2.Range("B$r").NumberFormat = $ws2.Range("B2").NumberFormat
}

$ws2.Range("A1:D17").Select() | Out-Null

# Leave "Schedule" as the active/visible tab with the author's final
# selection, matching the workbook's saved state.
$ws1.Activate()
$ws1.Range("F23").Select() | Out-Null
